$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B2").Value = 0
$ws.Range("F2").Value = 24.52000000000039
$ws.Range("H2").Value = 0.1420904863558349
$ws.Range("I2").Value = 0.1420904863558349
$ws.Range("L2").Value = 21.33312214223257
$ws.Range("M2").Value = "[-7.130940907468755, 49.79718519193389]"
$ws.Range("N2").Value = 0.1381557103380096
$ws.Range("O2").Value = 0.1381557103380096
$ws.Range("P2").Value = 1.163552834662887
$ws.Range("Q2").Value = "[-1.7044476659115784, 4.031553335237351]"
$ws.Range("R2").Value = 0.418162084492578
$ws.Range("S2").Value = 0.418162084492578
$ws.Range("T2").Value = 57.55063296840578
$ws.Range("U2").Value = "[42.433902379313324, 72.66736355749823]"
$ws.Range("V2").Value = [double]"1.05374664549629e-09"
$ws.Range("W2").Value = [double]"1.05374664549629e-09"
$ws.Range("X2").Value = 19.97925925925958
$ws.Range("Y2").Value = 8.78694694694709
$ws.Range("Z2").Value = 31.17157157157207
$ws.Range("F3").Value = 24.52000000000039
$ws.Range("H3").Value = 0.01587655987138681
$ws.Range("I3").Value = 0.01587655987138681
$ws.Range("L3").Value = 32.8804660640216
$ws.Range("M3").Value = "[5.582553122523343, 60.178379005519865]"
$ws.Range("N3").Value = 0.01933808008628191
$ws.Range("O3").Value = 0.01933808008628191
$ws.Range("P3").Value = 1.402552876377425
$ws.Range("Q3").Value = "[0.23271056693257552, 2.572395185822274]"
$ws.Range("R3").Value = 0.0198724409230302
$ws.Range("S3").Value = 0.0198724409230302
$ws.Range("T3").Value = 67.31151724737798
$ws.Range("U3").Value = "[51.648078086606844, 82.97495640814911]"
$ws.Range("V3").Value = [double]"3.911293511293934e-11"
$ws.Range("W3").Value = [double]"3.911293511293934e-11"
$ws.Range("X3").Value = 19.04656656656687
$ws.Range("Y3").Value = 14.48128128128151
$ws.Range("Z3").Value = 23.61185185185224
$ws.Range("F4").Value = 24.52000000000039
$ws.Range("H4").Value = 0.002794187975679763
$ws.Range("I4").Value = 0.002794187975679763
$ws.Range("L4").Value = 40.32315010676437
$ws.Range("M4").Value = "[12.398476117082268, 68.24782409644648]"
$ws.Range("N4").Value = 0.005625279566831898
$ws.Range("O4").Value = 0.005625279566831898
$ws.Range("P4").Value = 1.301921279866041
$ws.Range("Q4").Value = "[0.48428955821104047, 2.119553001521041]"
$ws.Range("R4").Value = 0.002470335975375937
$ws.Range("S4").Value = 0.002470335975375937
$ws.Range("T4").Value = 57.73680132759517
$ws.Range("U4").Value = "[41.97072764629502, 73.50287500889533]"
$ws.Range("V4").Value = [double]"2.833144607805593e-09"
$ws.Range("W4").Value = [double]"2.833144607805593e-09"
$ws.Range("X4").Value = 19.43927927927959
$ws.Range("Y4").Value = 16.24848848848875
$ws.Range("Z4").Value = 22.63007007007043
$ws.Range("F5").Value = 24.52000000000039
$ws.Range("H5").Value = [double]"2.427386996473047e-05"
$ws.Range("I5").Value = [double]"2.427386996473047e-05"
$ws.Range("L5").Value = 44.96067884505165
$ws.Range("M5").Value = "[21.531839632694457, 68.38951805740884]"
$ws.Range("N5").Value = 0.0003539273451895752
$ws.Range("O5").Value = 0.0003539273451895752
$ws.Range("P5").Value = 0.8993948938205012
$ws.Range("Q5").Value = "[0.37107901213573147, 1.427710775505271]"
$ws.Range("R5").Value = 0.001307671816263678
$ws.Range("S5").Value = 0.001307671816263678
$ws.Range("T5").Value = 46.8886813932255
$ws.Range("U5").Value = "[34.608893939775186, 59.16846884667581]"
$ws.Range("V5").Value = [double]"9.759397734399045e-10"
$ws.Range("W5").Value = [double]"9.759397734399045e-10"
$ws.Range("X5").Value = 21.01013013013047
$ws.Range("Y5").Value = 18.9483883883887
$ws.Range("Z5").Value = 23.07187187187224
$ws.Range("B6").Value = 1
$ws.Range("F6").Value = 25.8200000000006
$ws.Range("H6").Value = 0.0001375968508305903
$ws.Range("I6").Value = 0.0001375968508305903
$ws.Range("L6").Value = 50.82941859554889
$ws.Range("M6").Value = "[21.513492637322983, 80.14534455377479]"
$ws.Range("N6").Value = 0.001086255455417851
$ws.Range("O6").Value = 0.001086255455417851
$ws.Range("P6").Value = 0.8616580451287321
$ws.Range("Q6").Value = "[0.3333421634439624, 1.3899739268135018]"
$ws.Range("R6").Value = 0.001980527455344605
$ws.Range("S6").Value = 0.001980527455344605
$ws.Range("T6").Value = 64.9192575558852
$ws.Range("U6").Value = "[49.630043784130066, 80.20847132764032]"
$ws.Range("V6").Value = [double]"5.496425536932747e-11"
$ws.Range("W6").Value = [double]"5.496425536932747e-11"
$ws.Range("X6").Value = 22.27911911911963
$ws.Range("Y6").Value = 20.10806806806853
$ws.Range("Z6").Value = 24.45017017017074
$ws.Range("F7").Value = 25.8200000000006
$ws.Range("H7").Value = [double]"1.001502088748829e-05"
$ws.Range("I7").Value = [double]"1.001502088748829e-05"
$ws.Range("L7").Value = 63.51877427034402
$ws.Range("M7").Value = "[35.50730414799442, 91.53024439269362]"
$ws.Range("N7").Value = [double]"3.827302125070098e-05"
$ws.Range("O7").Value = [double]"3.827302125070098e-05"
$ws.Range("P7").Value = 0.1949737182408082
$ws.Range("Q7").Value = "[-0.3207632138800385, 0.7107106503616549]"
$ws.Range("R7").Value = 0.4503727851081925
$ws.Range("S7").Value = 0.4503727851081925
$ws.Range("T7").Value = 68.82611371979564
$ws.Range("U7").Value = "[52.949583062502114, 84.70264437708917]"
$ws.Range("V7").Value = [double]"3.04751779367507e-11"
$ws.Range("W7").Value = [double]"3.04751779367507e-11"
$ws.Range("X7").Value = 25.01877877877936
$ws.Range("Y7").Value = 22.89941941941995
$ws.Range("Z7").Value = 27.13813813813876
$ws.Range("F8").Value = 25.8200000000006
$ws.Range("H8").Value = 0.0004597407380754204
$ws.Range("I8").Value = 0.0004597407380754204
$ws.Range("L8").Value = 46.16002626225416
$ws.Range("M8").Value = "[21.97534904734583, 70.3447034771625]"
$ws.Range("N8").Value = 0.0003773712061734269
$ws.Range("O8").Value = 0.0003773712061734269
$ws.Range("P8").Value = 0.05660527303765317
$ws.Range("Q8").Value = "[-0.6100790538502725, 0.7232895999255788]"
$ws.Range("R8").Value = 0.8649834975903561
$ws.Range("S8").Value = 0.8649834975903561
$ws.Range("T8").Value = 68.38681649750841
$ws.Range("U8").Value = "[53.556231712075906, 83.2174012829409]"
$ws.Range("V8").Value = [double]"4.99977836909693e-12"
$ws.Range("W8").Value = [double]"4.99977836909693e-12"
$ws.Range("X8").Value = 25.58738738738798
$ws.Range("Y8").Value = 22.84772772772825
$ws.Range("Z8").Value = 28.32704704704771
$ws.Range("B9").Value = 0
$ws.Range("F9").Value = 25.8200000000006
$ws.Range("H9").Value = [double]"1.233102615016968e-06"
$ws.Range("I9").Value = [double]"1.233102615016968e-06"
$ws.Range("L9").Value = 54.84569182144914
$ws.Range("M9").Value = "[31.690861309377254, 78.00052233352103]"
$ws.Range("N9").Value = [double]"1.965790913560639e-05"
$ws.Range("O9").Value = [double]"1.965790913560639e-05"
$ws.Range("P9").Value = -0.4025263860455386
$ws.Range("Q9").Value = "[-0.8679475199106941, 0.0628947478196169]"
$ws.Range("R9").Value = 0.08835334470312461
$ws.Range("S9").Value = 0.08835334470312461
$ws.Range("T9").Value = 62.83015869490558
$ws.Range("U9").Value = "[50.115686215259565, 75.5446311745516]"
$ws.Range("V9").Value = [double]"6.0285110237146e-13"
$ws.Range("W9").Value = [double]"6.0285110237146e-13"
$ws.Range("X9").Value = 1.654134134134171
$ws.Range("Y9").Value = -0.2584584584584693
$ws.Range("Z9").Value = 3.566726726726812
$ws.Range("F10").Value = 25.8200000000006
$ws.Range("H10").Value = [double]"3.749629340354943e-06"
$ws.Range("I10").Value = [double]"3.749629340354943e-06"
$ws.Range("L10").Value = 68.2096059002437
$ws.Range("M10").Value = "[34.78126371094784, 101.63794808953956]"
$ws.Range("N10").Value = 0.0001655130496664992
$ws.Range("O10").Value = 0.0001655130496664992
$ws.Range("P10").Value = -1.031473864241693
$ws.Range("Q10").Value = "[-1.522052897234694, -0.5408948312486928]"
$ws.Range("R10").Value = 0.0001114828898529652
$ws.Range("S10").Value = 0.0001114828898529652
$ws.Range("T10").Value = 72.91864675502721
$ws.Range("U10").Value = "[55.5112378257745, 90.32605568427992]"
$ws.Range("V10").Value = [double]"8.040190735414399e-11"
$ws.Range("W10").Value = [double]"8.040190735414399e-11"
$ws.Range("X10").Value = 4.238718718718815
$ws.Range("Y10").Value = 2.222742742742792
$ws.Range("Z10").Value = 6.254694694694837
